# Refresh the cryptos price/volume snapshot (row 2..51, columns D=Price,
# E=Volume(1h)) to match the latest scrape, as produced by the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# Column D values that look numeric (e.g. "0.999") are entered with a
# leading apostrophe so Excel stores them as text (matching how every
# other cell in this text-only price column is stored) instead of
# silently converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '64.028.51'; E = '  +6.05%  ' },
    @{ Row = 3; D = '2.748.84'; E = '  +5.16%  ' },
    @{ Row = 4; D = '''0.999'; E = '  -0.05%  ' },
    @{ Row = 5; D = '''595.77'; E = '  +1.87%  ' },
    @{ Row = 6; D = '''153.23'; E = '  +6.92%  ' },
    @{ Row = 7; D = '''0.994'; E = '  -0.38%  ' },
    @{ Row = 8; D = '''0.613'; E = '  +2.88%  ' },
    @{ Row = 9; D = '2.792.02'; E = '  +6.44%  ' },
    @{ Row = 10; D = $null; E = '  +4.22%  ' },
    @{ Row = 11; D = $null; E = '  +8.88%  ' },
    @{ Row = 12; D = $null; E = '  +4.66%  ' },
    @{ Row = 13; D = $null; E = '  +1.74%  ' },
    @{ Row = 14; D = '3.228.83'; E = '  +4.98%  ' },
    @{ Row = 15; D = '''26.83'; E = '  +8.48%  ' },
    @{ Row = 16; D = '63.866.17'; E = '  +5.81%  ' },
    @{ Row = 17; D = $null; E = '  +9.24%  ' },
    @{ Row = 18; D = '2.765.79'; E = '  +5.82%  ' },
    @{ Row = 19; D = '''12.15'; E = '  +6.75%  ' },
    @{ Row = 20; D = '''4.93'; E = '  +6.27%  ' },
    @{ Row = 21; D = '''367.07'; E = '  +5.88%  ' },
    @{ Row = 22; D = $null; E = '  +2.00%  ' },
    @{ Row = 23; D = '''0.539'; E = '  +0.65%  ' },
    @{ Row = 24; D = '''0.996'; E = '  -0.26%  ' },
    @{ Row = 25; D = '''66.30'; E = '  +4.20%  ' },
    @{ Row = 26; D = $null; E = '  +5.49%  ' },
    @{ Row = 27; D = '''8.69'; E = '  +8.15%  ' },
    @{ Row = 28; D = '''0.998'; E = '  -0.08%  ' },
    @{ Row = 29; D = '0.0₃0905'; E = '  +13.58%  ' },
    @{ Row = 30; D = '''2.06'; E = '  +6.83%  ' },
    @{ Row = 31; D = '''7.18'; E = '  +11.38%  ' },
    @{ Row = 32; D = '''171.02'; E = '  +1.23%  ' },
    @{ Row = 33; D = '''1.20'; E = '  +18.33%  ' },
    @{ Row = 34; D = '''0.996'; E = '  -0.21%  ' },
    @{ Row = 35; D = '''20.75'; E = '  +6.48%  ' },
    @{ Row = 36; D = '''4.82'; E = '  +12.38%  ' },
    @{ Row = 37; D = $null; E = '  +11.64%  ' },
    @{ Row = 38; D = $null; E = '  +9.97%  ' },
    @{ Row = 39; D = '''1.03'; E = '  +21.04%  ' },
    @{ Row = 40; D = '''352.37'; E = '  +10.20%  ' },
    @{ Row = 41; D = $null; E = '  +9.13%  ' },
    @{ Row = 42; D = '''39.39'; E = '  +2.40%  ' },
    @{ Row = 43; D = '''5.69'; E = '  +12.08%  ' },
    @{ Row = 44; D = '''22.37'; E = '  +12.14%  ' },
    @{ Row = 45; D = '''144.45'; E = '  +6.14%  ' },
    @{ Row = 46; D = $null; E = '  +11.05%  ' },
    @{ Row = 47; D = '''0.0596'; E = '  +8.59%  ' },
    @{ Row = 48; D = '''0.653'; E = '  +7.42%  ' },
    @{ Row = 49; D = $null; E = '  +7.56%  ' },
    @{ Row = 50; D = $null; E = '  +2.89%  ' },
    @{ Row = 51; D = '2.181.24'; E = '  +7.70%  ' }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $ws.Range("D" + $u.Row).Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
